$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 36, pushing existing rows 36:111 down to 37:112
$ws.Rows.Item(36).Insert()

# Populate the new row 36 with the new data record
$ws.Range("A36").Value = 11
$ws.Range("B36").Value = "Vega Monumental Concepción"
$ws.Range("C36").Value = "Bíobío"
$ws.Range("D36").NumberFormat = $ws.Range("D37").NumberFormat
$ws.Range("D36").Value = 44498
$ws.Range("E36").Value = 8
$ws.Range("F36").Value = "Fruta"
$ws.Range("G36").Value = 100108
$ws.Range("H36").Value = "Tropicales y subtropicales"
$ws.Range("I36").Value = 100108005
$ws.Range("J36").Value = "Piña"
$ws.Range("K36").Value = "Caramelo"
$ws.Range("L36").Value = "Segunda"
$ws.Range("M36").Value = 200
$ws.Range("N36").Value = 19000
$ws.Range("O36").Value = 20000
$ws.Range("P36").Value = 19500
$ws.Range("Q36").Value = "$/caja 14 unidades"
$ws.Range("R36").Value = "Ecuador"
$ws.Range("S36").Value = 1393
$ws.Range("T36").Value = 14
